$d = $word.ActiveDocument

# 1) Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
    $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the mailing address line "2910 Lamory Pl, Santa Clara CA 95051" into
#    two separate paragraphs: "2910 Lamory Pl" and "Santa Clara, CA 95051".
#    Only the first (non-tabular) occurrence of this text should be touched.
$addrRange = $d.Content
$found = $addrRange.Find.Execute("2910 Lamory Pl, Santa Clara CA 95051")
if ($found) {
    $addrRange.Text = "2910 Lamory Pl`rSanta Clara, CA 95051"
}

# 3) Remove the now-redundant blank "No Spacing" paragraph that sits right
#    after the "Board of Directors" signature line.
$bodRange = $d.Content
$foundBod = $bodRange.Find.Execute("Board of Directors")
if ($foundBod) {
    $bodRange.Collapse(0) | Out-Null
    $bodRange.Move(1, 1) | Out-Null
    $bodRange.Expand(4) | Out-Null
    if ($bodRange.Text.Trim() -eq "") {
        $bodRange.Delete()
    }
}
